$wb = $excel.ActiveWorkbook

# 1. "Ready for handoff" -> "In Translation" (Status column, every locale sheet).
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        $v = $cell.Value2
        # Guard the type: PowerShell's -eq coerces a bare string RHS to
        # [bool] when $v is already a [bool] (TRUE/FALSE cells), which
        # would otherwise make every boolean cell match.
        if (($v -is [string]) -and ($v -eq "Ready for handoff")) {
            $cell.Value2 = "In Translation"
        }
    }
}

# 2. Narrow the "Status" columns from ~17.22 chars to ~13.41 chars.
#    (Overview!E:F and the per-locale sheets' column C.)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
